$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C29").Value = 103
$ws.Range("D29").Value = 19
$ws.Range("E29").Value = 84
$ws.Range("F29").Value = 3.270223752151463
